$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read current (pre-edit) values for the columns that are rotating among rows 2-5.
$rows = @(2, 3, 4, 5)
$idCol = "A"
$ostCol = "Q"
$nordCol = "R"
$startCol = "Y"
$endCol = "AA"

$idVals = @{}
$ostVals = @{}
$nordVals = @{}
$startVals = @{}
$endVals = @{}

foreach ($r in $rows) {
    $idVals[$r] = $ws.Range($idCol + $r).Value2
    $ostVals[$r] = $ws.Range($ostCol + $r).Value2
    $nordVals[$r] = $ws.Range($nordCol + $r).Value2
    $startVals[$r] = $ws.Range($startCol + $r).Text
    $endVals[$r] = $ws.Range($endCol + $r).Text
}

# New row values come from the row below (cyclic shift upward of data, i.e.
# row2<-row5(old), row3<-row2(old), row4<-row3(old), row5<-row4(old)).
$sourceFor = @{ 2 = 5; 3 = 2; 4 = 3; 5 = 4 }

foreach ($r in $rows) {
    $src = $sourceFor[$r]

    $ws.Range($idCol + $r).Value = $idVals[$src]
    $ws.Range($ostCol + $r).Value = $ostVals[$src]
    $ws.Range($nordCol + $r).Value = $nordVals[$src]

    # Date-like text must stay as text (not get auto-converted to a date
    # serial number with a date number format).
    $ws.Range($startCol + $r).NumberFormat = "@"
    $ws.Range($startCol + $r).Value = $startVals[$src]
    $ws.Range($startCol + $r).Style = "Normal"

    $ws.Range($endCol + $r).NumberFormat = "@"
    $ws.Range($endCol + $r).Value = $endVals[$src]
    $ws.Range($endCol + $r).Style = "Normal"
}
